$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="65.086.08"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)

$ws.Range("D3").Formula = '="3.179.00"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)

$ws.Range("E4").Formula = '="  -0.02%  "'
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)

$ws.Range("D5").Formula = '="579.13"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)

$ws.Range("E5").Formula = '="  +3.22%  "'
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)

$ws.Range("D6").Formula = '="151.36"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)

$ws.Range("E6").Formula = '="  +4.99%  "'
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)

$ws.Range("E7").Formula = '="  -0.05%  "'
$ws.Range("E7").Copy()
$ws.Range("E7").PasteSpecial(-4163)

$ws.Range("D8").Formula = '="3.178.46"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)

$ws.Range("E8").Formula = '="  +3.63%  "'
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)

$ws.Range("D9").Formula = '="0.531"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)

$ws.Range("E9").Formula = '="  +3.21%  "'
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial(-4163)

$ws.Range("E10").Formula = '="  +5.08%  "'
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)

$ws.Range("E11").Formula = '="  +0.76%  "'
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial(-4163)

$ws.Range("E12").Formula = '="  +2.02%  "'
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)

$ws.Range("E13").Formula = '="  +18.11%  "'
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)

$ws.Range("D14").Formula = '="38.11"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)

$ws.Range("E14").Formula = '="  +6.31%  "'
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)

$ws.Range("D15").Formula = '="3.698.90"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)

$ws.Range("E15").Formula = '="  +3.63%  "'
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)

$ws.Range("D16").Formula = '="65.153.52"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)

$ws.Range("E16").Formula = '="  +1.55%  "'
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)

$ws.Range("D17").Formula = '="3.173.47"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)

$ws.Range("E17").Formula = '="  +3.46%  "'
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial(-4163)

$ws.Range("D18").Formula = '="7.20"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)

$ws.Range("E18").Formula = '="  +5.62%  "'
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial(-4163)

$ws.Range("E19").Formula = '="  +1.35%  "'
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)

$ws.Range("D20").Formula = '="514.20"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)

$ws.Range("E21").Formula = '="  +6.71%  "'
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial(-4163)

$ws.Range("D22").Formula = '="0.735"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)

$ws.Range("D23").Formula = '="15.42"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)

$ws.Range("E23").Formula = '="  +6.76%  "'
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)

$ws.Range("D24").Formula = '="7.84"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)

$ws.Range("E24").Formula = '="  +3.35%  "'
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial(-4163)

$ws.Range("D25").Formula = '="85.49"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)

$ws.Range("E25").Formula = '="  +3.41%  "'
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial(-4163)

$ws.Range("D26").Formula = '="1.00"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)

$ws.Range("E26").Formula = '="  +0.04%  "'
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial(-4163)

$ws.Range("D27").Formula = '="9.04"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)

$ws.Range("E27").Formula = '="  +11.19%  "'
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial(-4163)

$ws.Range("E28").Formula = '="  +4.23%  "'
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial(-4163)

$ws.Range("E29").Formula = '="  +7.32%  "'
$ws.Range("E29").Copy()
$ws.Range("E29").PasteSpecial(-4163)

$ws.Range("D30").Formula = '="28.06"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)

$ws.Range("E30").Formula = '="  +6.48%  "'
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial(-4163)

$ws.Range("D31").Formula = '="2.80"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)

$ws.Range("E31").Formula = '="  +13.56%  "'
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial(-4163)

$ws.Range("D32").Formula = '="0.999"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)

$ws.Range("E32").Formula = '="  +0.02%  "'
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial(-4163)

$ws.Range("E33").Formula = '="  +5.49%  "'
$ws.Range("E33").Copy()
$ws.Range("E33").PasteSpecial(-4163)

$ws.Range("D34").Formula = '="6.26"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)

$ws.Range("E34").Formula = '="  +7.61%  "'
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial(-4163)

$ws.Range("E35").Formula = '="  +5.76%  "'
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial(-4163)

$ws.Range("D36").Formula = '="55.65"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)

$ws.Range("E36").Formula = '="  +1.28%  "'
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial(-4163)

$ws.Range("D37").Formula = '="0.0903"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)

$ws.Range("E37").Formula = '="  +10.40%  "'
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial(-4163)

$ws.Range("D38").Formula = '="478.37"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)

$ws.Range("E38").Formula = '="  +5.44%  "'
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial(-4163)

$ws.Range("D39").Formula = '="3.17"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)

$ws.Range("E39").Formula = '="  +11.90%  "'
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial(-4163)

$ws.Range("E40").Formula = '="  +2.22%  "'
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial(-4163)

$ws.Range("D41").Formula = '="8.67"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)

$ws.Range("E41").Formula = '="  +4.41%  "'
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial(-4163)

$ws.Range("D42").Formula = '="3.066.24"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)

$ws.Range("E42").Formula = '="  +1.22%  "'
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial(-4163)

$ws.Range("D43").Formula = '="0.120"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)

$ws.Range("E43").Formula = '="  +2.62%  "'
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial(-4163)

$ws.Range("E44").Formula = '="  +7.73%  "'
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial(-4163)

$ws.Range("E45").Formula = '="  +7.67%  "'
$ws.Range("E45").Copy()
$ws.Range("E45").PasteSpecial(-4163)

$ws.Range("D46").Formula = '="29.17"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)

$ws.Range("E46").Formula = '="  +4.80%  "'
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial(-4163)

$ws.Range("D47").Formula = '="0.0₃0611"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)

$ws.Range("E47").Formula = '="  +17.37%  "'
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial(-4163)

$ws.Range("E48").Formula = '="  -0.08%  "'
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial(-4163)

$ws.Range("E49").Formula = '="  +2.05%  "'
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial(-4163)

$ws.Range("E50").Formula = '="  +8.60%  "'
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial(-4163)

$ws.Range("D51").Formula = '="120.57"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)

$ws.Range("E51").Formula = '="  +1.54%  "'
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial(-4163)

$excel.CutCopyMode = 0
